# Refresh the cryptocurrency price table (Price / Volume(1h) columns) and
# correct the row 8/9 ordering (XRP <-> LidoStakedEther), per the commit's
# scraped-data update.
#
# Column D (Price) values are written with a leading apostrophe -- Excel's
# standard 'force text' entry marker -- so numeric-looking strings such as
# "1.00" stay text instead of being coerced to numbers; this matches the
# workbook's existing inline-string cells in that column. The apostrophe is
# prepended via string concatenation rather than embedded in a quoted
# literal, so it can never be mistaken for a quote-escape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '62.687.87'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = "'" + '3.025.22'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'" + '585.53'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").Value = "'" + '147.79'
$ws.Range("E6").Value = '  -4.43%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = "'" + '3.025.73'
$ws.Range("E8").Value = '  -1.17%  '
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = "'" + '0.522'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("D11").Value = "'" + '5.78'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("D14").Value = "'" + '34.92'
$ws.Range("E14").Value = '  -4.99%  '
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").Value = "'" + '3.524.55'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = "'" + '62.665.06'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = "'" + '3.024.07'
$ws.Range("E19").Value = '  -1.29%  '
$ws.Range("D20").Value = "'" + '465.25'
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Value = "'" + '13.99'
$ws.Range("E21").Value = '  -1.91%  '
$ws.Range("E22").Value = '  -1.86%  '
$ws.Range("D23").Value = "'" + '7.39'
$ws.Range("E23").Value = '  -1.06%  '
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("D25").Value = "'" + '80.19'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").Value = "'" + '12.41'
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("D27").Value = "'" + '10.17'
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").Value = "'" + '1.00'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").Value = "'" + '7.18'
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").Value = "'" + '27.51'
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("D35").Value = "'" + '1.03'
$ws.Range("E35").Value = '  -0.70%  '
$ws.Range("D36").Value = "'" + '0.0₃0797'
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("E38").Value = '  -2.45%  '
$ws.Range("D39").Value = "'" + '50.53'
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = "'" + '9.03'
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("E41").Value = '  -8.70%  '
$ws.Range("D42").Value = "'" + '424.68'
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").Value = "'" + '2.788.97'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = "'" + '0.0354'
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("D47").Value = "'" + '37.99'
$ws.Range("E47").Value = '  -6.48%  '
$ws.Range("D48").Value = "'" + '129.71'
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = "'" + '24.20'
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("E51").Value = '  -0.47%  '
